$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J4 ("View Mail" status) changes from "In-Progress" to "Done".
# (The shared-string table is compacted automatically on save, so the
# now-unused "In-Progress" entry disappears and K4/L4 - which were and
# remain "TODO" - are re-pointed to the shifted index.)
$ws.Range("J4").Value = "Done"

# The sheet's active selection moved from J10 to J9.
$ws.Range("J9").Select()
